# Update race-result stats on the "Drivers" sheet, then make it the
# active/displayed sheet (matching the author's last view state).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Drivers")

# --- Updated driver stat values (Acceleration column, mostly) ---
$ws.Range("D13").Value = 81
$ws.Range("E15").Value = 78
$ws.Range("E16").Value = 74
$ws.Range("D19").Value = 74

# --- View state: Drivers becomes the active/selected sheet, zoomed to 85%,
#     with D14 selected (previously "Race Time Calculation" was active) ---
$ws.Activate()
$excel.ActiveWindow.Zoom = 85
$ws.Range("D14").Select()
